$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Metadata sheet: bump the publication Date value.
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(8, 2).Value = "2026-02-05T08:09:31+00:00"

# ---------------------------------------------------------------------------
# 2) Elements sheet: insert a new "fr-lm-evaluation.valeur" row right after
#    the "date" row (row 8) and before the "interpretation" row, shifting
#    every following row down by one.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Elements")

# Insert a fresh row at position 8 (old row 8 "interpretation" and everything
# below slides down to row 9+).
$ws.Rows.Item(8).Insert()

# Copy the formatting (borders / fill / wrap alignment) from the row that
# used to be directly above (row 7, "date") onto the newly inserted row so
# that its style matches every other data row instead of Excel's bare
# insert-row default.
$ws.Range("A7:AJ7").Copy()
$ws.Range("A8:AJ8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new row with the "Valeur de l'évaluation" element.
$ws.Cells.Item(8, 1).Value  = "fr-lm-evaluation.valeur"        # A - ID
$ws.Cells.Item(8, 2).Value  = "fr-lm-evaluation.valeur"        # B - Path
$ws.Cells.Item(8, 4).Value  = ""                                # D - Alias(s)
$ws.Cells.Item(8, 6).Value  = "1"                               # F - Min
$ws.Cells.Item(8, 7).Value  = "1"                               # G - Max
$ws.Cells.Item(8, 8).Value  = ""                                # H - Must Support?
$ws.Cells.Item(8, 9).Value  = ""                                # I - Is Modifier?
$ws.Cells.Item(8, 10).Value = ""                                # J - Is Summary?
$ws.Cells.Item(8, 11).Value = "Base`n"                          # K - Type(s)
$ws.Cells.Item(8, 12).Value = "Valeur de l'évaluation"          # L - Short
$ws.Cells.Item(8, 13).Value = "Valeur de l'évaluation"          # M - Definition
$ws.Cells.Item(8, 16).Value = ""                                # P - Default Value
$ws.Cells.Item(8, 18).Value = ""                                # R - Fixed Value
$ws.Cells.Item(8, 19).Value = ""                                # S - Pattern
$ws.Cells.Item(8, 20).Value = ""                                # T - Example
$ws.Cells.Item(8, 21).Value = ""                                # U - Minimum Value
$ws.Cells.Item(8, 22).Value = ""                                # V - Maximum Value
$ws.Cells.Item(8, 23).Value = ""                                # W - Maximum Length
$ws.Cells.Item(8, 24).Value = ""                                # X - Binding Strength
$ws.Cells.Item(8, 25).Value = ""                                # Y - Binding Description
$ws.Cells.Item(8, 26).Value = ""                                # Z - Binding Value Set
$ws.Cells.Item(8, 27).Value = ""                                # AA - Code
$ws.Cells.Item(8, 28).Value = ""                                # AB - Slicing Discriminator
$ws.Cells.Item(8, 29).Value = ""                                # AC - Slicing Description
$ws.Cells.Item(8, 30).Value = ""                                # AD - Slicing Ordered
$ws.Cells.Item(8, 31).Value = ""                                # AE - Slicing Rules
$ws.Cells.Item(8, 32).Value = "fr-lm-evaluation.valeur"         # AF - Base Path
$ws.Cells.Item(8, 33).Value = "1"                               # AG - Base Min
$ws.Cells.Item(8, 34).Value = "1"                               # AH - Base Max
$ws.Cells.Item(8, 35).Value = ""                                # AI - Condition(s)
$ws.Cells.Item(8, 36).Value = ""                                # AJ - Constraint(s)
